$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Calculations" worksheet between "About" and "CApULAbIFM".
# ---------------------------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")
$calcSheet  = $wb.Worksheets.Add([System.Type]::Missing, $aboutSheet)
$calcSheet.Name = "Calculations"

$calcSheet.Range("A1").Value = "Increased Annual CO2 Sequestration Achievable by Improved Management Practices per Acre"
$calcSheet.Range("A1").Font.Bold = $true

$calcSheet.Range("A2").Value = 2.1
$calcSheet.Range("B2").Value = "tons CO2 / acre / yr"
$calcSheet.Range("C2").Value = "Low Estimate"

$calcSheet.Range("A3").Value = 3.1
$calcSheet.Range("B3").Value = "tons CO2 / acre / yr"
$calcSheet.Range("C3").Value = "High Estimate"

$calcSheet.Range("A4").Formula = "=AVERAGE(A2:A3)"
$calcSheet.Range("B4").Value = "tons CO2 / acre / yr"
$calcSheet.Range("C4").Value = "Average"

$calcSheet.Range("A6").Formula = "=A4*10^6"
$calcSheet.Range("A6").NumberFormat = "0.00E+00"
$calcSheet.Range("B6").Value = "g CO2 / acre / yr"
$calcSheet.Range("C6").Value = "Average, converted to grams CO2"

$calcSheet.Range("B1").ColumnWidth = 18.1666666

# ---------------------------------------------------------------------------
# 2. Update the "About" sheet: change the Source text and add citation rows.
#    The old "Notes" / "We do not use..." remark (rows 5-6) is no longer
#    applicable, so new rows are inserted above it (pushing it down) and its
#    text is then cleared out, leaving the citation block in its place.
# ---------------------------------------------------------------------------
$aboutSheet.Range("B3").Value = "U.S. EPA"

$aboutSheet.Range("A4:A7").EntireRow.Insert()

$aboutSheet.Range("B4").Value = 2005
$aboutSheet.Range("B4").HorizontalAlignment = -4131

$aboutSheet.Range("B5").Value = "Greenhouse Gas Mitigation Potential in U.S. Forestry and Agriculture"

$aboutSheet.Range("B6").Value = "http://www.epa.gov/climate/climatechange/Downloads/ccs/ghg_mitigation_forestry_ag_2005.pdf"
$aboutSheet.Hyperlinks.Add($aboutSheet.Range("B6"), "http://www.epa.gov/climate/climatechange/Downloads/ccs/ghg_mitigation_forestry_ag_2005.pdf")

$aboutSheet.Range("B7").Value = "Page 2-3, Table 2-1"

# Rows 5-6 ("Notes" / explanatory text) were pushed down to rows 9-10 by the
# insert above; clear their now-obsolete contents.
$aboutSheet.Range("A9").ClearContents()
$aboutSheet.Range("A10").ClearContents()
$aboutSheet.Range("A10").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Update the "CApULAbIFM" sheet: point the abated-CO2 cell at the new
#    Calculations sheet instead of the old hard-coded zero / "not used" note.
#    (B1 "Per Acre" / A2 "CO2 Abated (g)" text itself is unchanged - only
#    B2's value becomes a live formula.)
# ---------------------------------------------------------------------------
$mainSheet = $wb.Worksheets.Item("CApULAbIFM")
$mainSheet.Range("B2").Formula = "=Calculations!A6"
$mainSheet.Range("B2").NumberFormat = "0.00E+00"
